$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ccl21b", "Ackr4", "ECs", 2, 0.6666666666666666, 0.07908133333333334, 0.237244, 0.1380838080781507, 0.1380838080781507, 1, 0.3333333333333333, 0.377371, 1.132113, 0.4698794580655765, 0.4698794580655764, 0.02984300184133333, 0.268587016572, 0.06488274490739254, 0.06488274490739251),
    @("ECs", "Ccl21b", "Ackr4", "FAPs", 2, 0.6666666666666666, 0.07908133333333334, 0.237244, 0.1380838080781507, 0.1380838080781507, 1, 0.3333333333333333, 0.3560133333333333, 1.06804, 0.443286188209444, 0.443286188209444, 0.02815400908444445, 0.25338608176, 0.06121064493640788, 0.06121064493640786),
    @("ECs", "Ccl21b", "Ackr4", "MuSCs", 2, 0.6666666666666666, 0.07908133333333334, 0.237244, 0.1380838080781507, 0.1380838080781507, 2, 0.6666666666666666, 0.06973866666666667, 0.209216, 0.08683435372497944, 0.08683435372497944, 0.005515026744888889, 0.04963524070400001, 0.01199041823435031, 0.01199041823435031),
    @("FAPs", "Ccl21b", "Ackr4", "ECs", 2, 0.6666666666666666, 0.07720866666666666, 0.231626, 0.1348139473702591, 0.134813947370259, 1, 0.3333333333333333, 0.377371, 1.132113, 0.4698794580655765, 0.4698794580655764, 0.02913631174866666, 0.262226805738, 0.06334630453001848, 0.06334630453001845),
    @("FAPs", "Ccl21b", "Ackr4", "FAPs", 2, 0.6666666666666666, 0.07720866666666666, 0.231626, 0.1348139473702591, 0.134813947370259, 1, 0.3333333333333333, 0.3560133333333333, 1.06804, 0.443286188209444, 0.443286188209444, 0.02748731478222222, 0.24738583304, 0.05976116084723074, 0.05976116084723073),
    @("FAPs", "Ccl21b", "Ackr4", "MuSCs", 2, 0.6666666666666666, 0.07720866666666666, 0.231626, 0.1348139473702591, 0.134813947370259, 2, 0.6666666666666666, 0.06973866666666667, 0.209216, 0.08683435372497944, 0.08683435372497944, 0.005384429468444444, 0.048459865216, 0.01170648199300984, 0.01170648199300984),
    @("Inflammatory-Mac", "Ccl21b", "Ackr4", "ECs", 1, 0.3333333333333333, 0.03917266666666667, 0.117518, 0.06839933974190335, 0.06839933974190333, 1, 0.3333333333333333, 0.377371, 1.132113, 0.4698794580655765, 0.4698794580655764, 0.01478262839266667, 0.133043655534, 0.03213944468996879, 0.03213944468996878),
    @("Inflammatory-Mac", "Ccl21b", "Ackr4", "FAPs", 1, 0.3333333333333333, 0.03917266666666667, 0.117518, 0.06839933974190335, 0.06839933974190333, 1, 0.3333333333333333, 0.3560133333333333, 1.06804, 0.443286188209444, 0.443286188209444, 0.01394599163555556, 0.12551392472, 0.03032048259023107, 0.03032048259023107),
    @("Inflammatory-Mac", "Ccl21b", "Ackr4", "MuSCs", 1, 0.3333333333333333, 0.03917266666666667, 0.117518, 0.06839933974190335, 0.06839933974190333, 2, 0.6666666666666666, 0.06973866666666667, 0.209216, 0.08683435372497944, 0.08683435372497944, 0.002731849543111111, 0.024586645888, 0.005939412461703479, 0.005939412461703478),
    @("MuSCs", "Ccl21b", "Ackr4", "ECs", 3, 1, 0.3772426666666667, 1.131728, 0.658702904809687, 0.6587029048096869, 1, 0.3333333333333333, 0.377371, 1.132113, 0.4698794580655765, 0.4698794580655764, 0.1423604423626667, 1.281243981264, 0.3095109639381967, 0.3095109639381966),
    @("MuSCs", "Ccl21b", "Ackr4", "FAPs", 3, 1, 0.3772426666666667, 1.131728, 0.658702904809687, 0.6587029048096869, 1, 0.3333333333333333, 0.3560133333333333, 1.06804, 0.443286188209444, 0.443286188209444, 0.1343034192355556, 1.20873077312, 0.2919938998355744, 0.2919938998355743),
    @("MuSCs", "Ccl21b", "Ackr4", "MuSCs", 3, 1, 0.3772426666666667, 1.131728, 0.658702904809687, 0.6587029048096869, 2, 0.6666666666666666, 0.06973866666666667, 0.209216, 0.08683435372497944, 0.08683435372497944, 0.02630840058311111, 0.236775605248, 0.05719804103591582, 0.05719804103591581),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$wb.Save()